$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ r=60; A='November08  12:34:04'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=1; O=0.06318670014540355; P=0.08333185275395712; Q=0.2380952380952381; R=0.32; S=2.167353371031276; T=1.70749979664876; U=1; V=0.2380952380952381; W=0.32; X=2.167353371031276; Y=1.70749979664876 },
  @{ r=61; A='November08  12:58:37'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=1; O=0.06463994381446687; P=0.07156804349687364; Q=0.248015873015873; R=0.2844444444444444; S=2.188760974426994; T=1.742284579382931; U=1; V=0.248015873015873; W=0.2844444444444444; X=2.188760974426994; Y=1.742284579382931 },
  @{ r=62; A='November08  13:00:09'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=3; O=0.05018508398816699; P=0.07594142463472155; Q=0.4474206349206349; R=0.3377777777777778; S=1.569172261831564; T=2.092313976863362; U=3; V=0.4474206349206349; W=0.3377777777777778; X=1.569172261831564; Y=2.092313976863362 },
  @{ r=63; A='November08  13:02:49'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=1; O=0.06712434795640763; P=0.1484774870342679; Q=0.2361111111111111; R=0.2088888888888889; S=2.342033101307291; T=2.318045153428494; U=1; V=0.2361111111111111; W=0.2088888888888889; X=2.342033101307291; Y=2.318045153428494 },
  @{ r=64; A='November08  13:03:49'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; empty=$true },
  @{ r=65; A='November08  13:04:11'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=1; O=0.06717570610935726; P=0.08469244215223524; Q=0.244047619047619; R=0.2888888888888889; S=2.251542681017716; T=1.967513941783161; U=1; V=0.244047619047619; W=0.2888888888888889; X=2.251542681017716; Y=1.967513941783161 },
  @{ r=66; A='November08  13:04:51'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.06184496233860651; P=0.1136452208624946; Q=0.3571428571428572; R=0.2844444444444444; S=1.827914068349451; T=2.077391526784384; U=2; V=0.3571428571428572; W=0.2844444444444444; X=1.827914068349451; Y=2.077391526784384 },
  @{ r=67; A='November08  16:35:36'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; empty=$true },
  @{ r=68; A='November08  16:36:05'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; empty=$true },
  @{ r=69; A='November08  16:36:43'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; empty=$true },
  @{ r=70; A='November08  18:07:07'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.05952280298584983; P=0.1204833539326986; Q=0.3690476190476191; R=0.2755555555555556; S=1.716516405813988; T=1.956186766816161; U=2; V=0.3690476190476191; W=0.2755555555555556; X=1.716516405813988; Y=1.956186766816161 },
  @{ r=71; A='November08  18:16:02'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=38; O=0.003310799798262971; P=0.1120523161358303; Q=0.9742063492063492; R=0.3866666666666667; S=0.3362963545672746; T=1.50406855635713; U=38; V=0.9742063492063492; W=0.3866666666666667; X=0.3362963545672746; Y=1.50406855635713 },
  @{ r=72; A='November08  18:45:25'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.05038585904098693; P=0.08023414505852593; Q=0.3928571428571428; R=0.2933333333333333; S=1.56378915072861; T=1.717879830230017; U=2; V=0.3928571428571428; W=0.2933333333333333; X=1.56378915072861; Y=1.717879830230017 },
  @{ r=73; A='November08  18:47:10'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.05471803957507724; P=0.09047404792573717; Q=0.3849206349206349; R=0.28; S=1.786190412715338; T=1.907878402833891; U=2; V=0.3849206349206349; W=0.28; X=1.786190412715338; Y=1.907878402833891 },
  @{ r=74; A='November08  18:48:33'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.05069899298841991; P=0.0839737515979343; Q=0.4236111111111111; R=0.3377777777777778; S=1.567274451507842; T=1.952775802116908; U=2; V=0.4236111111111111; W=0.3377777777777778; X=1.567274451507842; Y=1.952775802116908 },
  @{ r=75; A='November08  18:50:35'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.05583590564746705; P=0.08431223842832777; Q=0.3998015873015873; R=0.2977777777777778; S=1.654119436773267; T=1.916014381759989; U=2; V=0.3998015873015873; W=0.2977777777777778; X=1.654119436773267; Y=1.916014381759989 },
  @{ r=76; A='November08  18:51:13'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f02f072fae8>'; K='10'; L='True'; M='32'; N=2; O=0.0602247439443119; P=0.09726434787114462; Q=0.3809523809523809; R=0.2711111111111111; S=1.630257064158317; T=2.211083193570267; U=2; V=0.3809523809523809; W=0.2711111111111111; X=1.630257064158317; Y=2.211083193570267 },
  @{ r=77; A='November08  19:15:17'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f974642ca60>'; K='10'; L='True'; M='32'; empty=$true },
  @{ r=78; A='November08  19:15:47'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f974642ca60>'; K='10'; L='True'; M='32'; N=2; O=0.03837354020770739; P=0.01966428215292925; Q=0.300081103000811; R=0.08110300081103; S=1.362217551845654; T=0.9285068112950217; U=2; V=0.300081103000811; W=0.08110300081103; X=1.362217551845654; Y=0.9285068112950217 },
  @{ r=79; A='November08  19:16:08'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f974642ca60>'; K='10'; L='True'; M='32'; N=1; O=0.04628054416962784; P=0.01478456822431677; Q=0.2165450121654501; R=0.08921330089213302; S=1.799792724841541; T=0.8075073538492893; U=1; V=0.2165450121654501; W=0.08921330089213302; X=1.799792724841541; Y=0.8075073538492893 },
  @{ r=80; A='November08  19:16:29'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f974642ca60>'; K='10'; L='True'; M='32'; N=2; O=0.03909223029406594; P=0.01627698277228263; Q=0.2984590429845904; R=0.06974858069748581; S=1.364596972963537; T=0.9509460862966862; U=2; V=0.2984590429845904; W=0.06974858069748581; X=1.364596972963537; Y=0.9509460862966862 },
  @{ r=81; A='November08  19:16:50'; B='resnet18'; C='200'; D='1'; E='True'; F='sgd'; G='1.0xsingle + 0.0Xmulti'; H='0.01'; I='0.9'; J='<function exp_lr_scheduler at 0x7f974642ca60>'; K='10'; L='True'; M='32'; N=1; O=0.048717105959067; P=0.01589492405708109; Q=0.1857258718572587; R=0.064882400648824; S=1.912735310121493; T=0.8717425753049757; U=1; V=0.1857258718572587; W=0.064882400648824; X=1.912735310121493; Y=0.8717425753049757 }
)

# First pass: set the text (shared-string) columns, forcing text storage
# by temporarily switching the number format to Text ("@") so that
# numeric-looking strings like "200" or "0.01" are stored as shared
# strings rather than being coerced to numbers.
$textCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")

foreach ($row in $rows) {
    $r = $row.r
    foreach ($col in $textCols) {
        $cell = $ws.Range($col + $r)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col]
        $cell.NumberFormat = "General"
    }

    if ($row.ContainsKey("empty")) {
        $sCell = $ws.Range("S" + $r)
        $sCell.NumberFormat = "@"
        $sCell.Value = " "
        $sCell.NumberFormat = "General"
    } else {
        $numCols = @("N","O","P","Q","R","S","T","U","V","W","X","Y")
        foreach ($col in $numCols) {
            $ws.Range($col + $r).Value = $row[$col]
        }
    }
}
